# Daily attendance processing - 2025-12-16 23:50:00
# Normalize the "Recorded By" (column G) attribution order for each
# session row: the list of recorders is re-ordered (reversed) for the
# affected rows to reflect the latest processing pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'G2' = 'backup@backdoor.com, System, system'
    'G3' = 'dnasr281@gmail.com, System'
    'G5' = 'backup@backdoor.com, System'
    'G6' = 'dnasr281@gmail.com, System'
    'G8' = 'backup@backdoor.com, System'
    'G10' = 'dnasr281@gmail.com, System'
    'G12' = 'dnasr281@gmail.com, System'
    'G13' = 'dnasr281@gmail.com, System'
    'G14' = 'dnasr281@gmail.com, System'
    'G15' = 'dnasr281@gmail.com, System'
    'G18' = 'dnasr281@gmail.com, System'
    'G19' = 'dnasr281@gmail.com, System'
    'G20' = 'dnasr281@gmail.com, System'
    'G21' = 'dnasr281@gmail.com, System'
    'G22' = 'dnasr281@gmail.com, System'
    'G24' = 'dnasr281@gmail.com, System'
    'G26' = 'dnasr281@gmail.com, System'
    'G28' = 'backup@backdoor.com, System, system'
    'G29' = 'dnasr281@gmail.com, System'
    'G31' = 'backup@backdoor.com, System'
    'G32' = 'dnasr281@gmail.com, System'
    'G34' = 'backup@backdoor.com, System'
    'G36' = 'dnasr281@gmail.com, System'
    'G38' = 'dnasr281@gmail.com, System'
    'G39' = 'dnasr281@gmail.com, System'
    'G40' = 'dnasr281@gmail.com, System'
    'G41' = 'dnasr281@gmail.com, System'
    'G44' = 'dnasr281@gmail.com, System'
    'G45' = 'dnasr281@gmail.com, System'
    'G46' = 'dnasr281@gmail.com, System'
    'G47' = 'dnasr281@gmail.com, System'
    'G48' = 'dnasr281@gmail.com, System'
    'G50' = 'dnasr281@gmail.com, System'
    'G52' = 'dnasr281@gmail.com, System'
    'G54' = 'backup@backdoor.com, System, system'
    'G55' = 'dnasr281@gmail.com, System'
    'G57' = 'backup@backdoor.com, System'
    'G58' = 'dnasr281@gmail.com, System'
    'G60' = 'backup@backdoor.com, System'
    'G62' = 'dnasr281@gmail.com, System'
    'G64' = 'dnasr281@gmail.com, System'
    'G65' = 'dnasr281@gmail.com, System'
    'G66' = 'dnasr281@gmail.com, System'
    'G67' = 'dnasr281@gmail.com, System'
    'G70' = 'dnasr281@gmail.com, System'
    'G71' = 'dnasr281@gmail.com, System'
    'G72' = 'dnasr281@gmail.com, System'
    'G73' = 'dnasr281@gmail.com, System'
    'G74' = 'dnasr281@gmail.com, System'
    'G76' = 'dnasr281@gmail.com, System'
    'G78' = 'dnasr281@gmail.com, System'
    'G80' = 'backup@backdoor.com, System'
    'G81' = 'backup@backdoor.com, System'
    'G82' = 'backup@backdoor.com, System'
    'G83' = 'dnasr281@gmail.com, System'
    'G84' = 'dnasr281@gmail.com, System'
    'G85' = 'dnasr281@gmail.com, System'
    'G86' = 'dnasr281@gmail.com, System'
    'G87' = 'dnasr281@gmail.com, admin@admin.com'
    'G90' = 'dnasr281@gmail.com, System'
    'G92' = 'dnasr281@gmail.com, System'
    'G99' = 'dnasr281@gmail.com, System'
    'G101' = 'dnasr281@gmail.com, System'
    'G106' = 'backup@backdoor.com, System'
    'G107' = 'backup@backdoor.com, System'
    'G108' = 'backup@backdoor.com, System'
    'G109' = 'dnasr281@gmail.com, System'
    'G110' = 'dnasr281@gmail.com, System'
    'G111' = 'dnasr281@gmail.com, System'
    'G112' = 'dnasr281@gmail.com, System'
    'G113' = 'dnasr281@gmail.com, admin@admin.com'
    'G116' = 'dnasr281@gmail.com, System'
    'G118' = 'dnasr281@gmail.com, System'
    'G125' = 'dnasr281@gmail.com, System'
    'G127' = 'dnasr281@gmail.com, System'
    'G132' = 'backup@backdoor.com, System'
    'G133' = 'backup@backdoor.com, System'
    'G134' = 'backup@backdoor.com, System'
    'G135' = 'dnasr281@gmail.com, System'
    'G136' = 'dnasr281@gmail.com, System'
    'G137' = 'dnasr281@gmail.com, System'
    'G138' = 'dnasr281@gmail.com, System'
    'G139' = 'dnasr281@gmail.com, admin@admin.com'
    'G142' = 'dnasr281@gmail.com, System'
    'G144' = 'dnasr281@gmail.com, System'
    'G151' = 'dnasr281@gmail.com, System'
    'G153' = 'dnasr281@gmail.com, System'
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
